$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: sequence numbers 1-5 in F1:J1
$ws.Range("F1").Value = 1
$ws.Range("G1").Value = 2
$ws.Range("H1").Value = 3
$ws.Range("I1").Value = 4
$ws.Range("J1").Value = 5

# Row 2: COUNTIF(.., "Cold Sore") across columns B..F, summed in L2
$ws.Range("F2").Formula = "=COUNTIF(B:B, ""Cold Sore"")"
$ws.Range("G2:J2").Formula = "=COUNTIF(C:C, ""Cold Sore"")"
$ws.Range("L2").Formula = "=SUM(F2:J2)"

# Row 3: COUNTIF(.., "*") across columns B..F, summed in L3
$ws.Range("F3").Formula = "=COUNTIF(B:B, ""*"")"
$ws.Range("G3:J3").Formula = "=COUNTIF(C:C, ""*"")"
$ws.Range("L3").Formula = "=SUM(F3:J3)"

# Row 4: ratio row2/row3, percentage formatted, including L4 (K4 stays empty)
$ws.Range("F4").Formula = "=F2/F3"
$ws.Range("G4:J4").Formula = "=G2/G3"
$ws.Range("L4").Formula = "=L2/L3"
$ws.Range("F4:L4").NumberFormat = "0.00%"
$ws.Range("K4").ClearContents()

# Match final view state (zoom + selection) seen in the authored workbook
$excel.ActiveWindow.Zoom = 144
[void]$ws.Range("K7").Select()
